$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D6: mark "OK" with the same style as D2/D4/D5 (centered)
$ws.Range("D6").Value = "OK"
$ws.Range("D6").HorizontalAlignment = -4108  # xlCenter

# Row 11: new todo item #10 - bullet shatter effect optimization
$ws.Range("A11").Value = 10
$ws.Range("A11").HorizontalAlignment = -4108  # xlCenter

$ws.Range("B6").Copy()
$ws.Range("B11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B11").Value = "子弹碎裂特效优化：衰减的子弹相应弱化碎裂效果"

# Move selection to match the recorded cursor position after the edit
$ws.Range("C22").Select()
